$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("C2").Value = 11.0010460251046
$ws.Range("D2").Value = 10043.0656934307
$ws.Range("F2").Value = 0.904472866026003
$ws.Range("I2").Value = 44
$ws.Range("J2").Value = 244

# Row 3
$ws.Range("C3").Value = 14.6908127208481
$ws.Range("D3").Value = 4219.81132075472
$ws.Range("E3").Value = 0.0909090909090909
$ws.Range("F3").Value = 3.18233096085409
$ws.Range("G3").Value = 265
$ws.Range("H3").Value = 21
$ws.Range("I3").Value = 44
$ws.Range("J3").Value = 583

# Row 4
$ws.Range("C4").Value = 22.135757717493
$ws.Range("D4").Value = 2151.3093289689
$ws.Range("E4").Value = 0.209345813617177
$ws.Range("F4").Value = 14.691084452975
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 53
$ws.Range("I4").Value = 44
$ws.Range("J4").Value = 1095

# Row 5
$ws.Range("C5").Value = 31.2593312597201
$ws.Range("D5").Value = 1287.77215189873
$ws.Range("E5").Value = 13.9098438652569
$ws.Range("F5").Value = 27.4744821731749
$ws.Range("G5").Value = 269
$ws.Range("H5").Value = 75
$ws.Range("I5").Value = 44
$ws.Range("J5").Value = 664

# Row 6
$ws.Range("C6").Value = 14.1897123893805
$ws.Range("D6").Value = 13206.9095477387
$ws.Range("F6").Value = 0.336244541484716
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 50
$ws.Range("J6").Value = 231

# Row 7
$ws.Range("C7").Value = 19.4140763765542
$ws.Range("D7").Value = 4762.58426966292
$ws.Range("E7").Value = 0.0388642252711808
$ws.Range("F7").Value = 5.40421792618629
$ws.Range("G7").Value = 139
$ws.Range("H7").Value = 15
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 584

# Row 8
$ws.Range("C8").Value = 25.8899572649573
$ws.Range("D8").Value = 2294.97734138973
$ws.Range("E8").Value = 0.624134969207516
$ws.Range("F8").Value = 23.1418269230769
$ws.Range("G8").Value = 189
$ws.Range("H8").Value = 19
$ws.Range("I8").Value = 50
$ws.Range("J8").Value = 851

# Row 9
$ws.Range("C9").Value = 37.0482261640798
$ws.Range("D9").Value = 1294.17302798982
$ws.Range("E9").Value = 14.9465836291888
$ws.Range("F9").Value = 45.6902654867257
$ws.Range("G9").Value = 69
$ws.Range("H9").Value = 10
$ws.Range("I9").Value = 50
$ws.Range("J9").Value = 462

# Row 12
$ws.Range("E12").Value = 2.24121951219512

# Row 14
$ws.Range("F14").Value = 0.0666666666666667
$ws.Range("H14").Value = 19

# Row 15
$ws.Range("E15").Value = 0.0132334581772784
$ws.Range("F15").Value = 0.645502645502645
$ws.Range("H15").Value = 78

# Row 16
$ws.Range("E16").Value = 0.187749287749288
$ws.Range("F16").Value = 5.74691358024691
$ws.Range("H16").Value = 108

# Row 17
$ws.Range("E17").Value = 6.19177793816913
$ws.Range("F17").Value = 27.3851351351351
$ws.Range("H17").Value = 15
